$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Orden de tipeo original: columna B (36, 39), columna C (36..39 salteando 36 ya hecha),
# luego columna A (36..39), y por ultimo D/E/F/G fila por fila.
$ws.Range("B36").Value = "Ver listado de blogs (user)"
$ws.Range("B39").Value = "Eliminar blog"
$ws.Range("C36").Value = "Ver listado de blogs desde el link de home"
$ws.Range("B37").Value = "Ver detalle de blog (admin)"
$ws.Range("B38").Value = "Ver detalle de blog (user)"
$ws.Range("C37").Value = "Ver delalle de blog desde el menu propio de admin"
$ws.Range("C38").Value = "Ver delalle de blog desde el listado que ve el user"
$ws.Range("C39").Value = "Eliminar Blog desde el menu propio de admin"
$ws.Range("A36").Value = "Caso #30"
$ws.Range("A37").Value = "Caso #31"
$ws.Range("A38").Value = "Caso #32"
$ws.Range("A39").Value = "Caso #33"

# Fila 36
$ws.Range("D36").Value = "12/17/2023"
$ws.Range("E36").Value = "SI"
$ws.Range("F36").Value = "-"
$ws.Range("G36").Value = "OK"

# Fila 37
$ws.Range("D37").Value = "12/17/2023"
$ws.Range("E37").Value = "SI"
$ws.Range("F37").Value = "-"
$ws.Range("G37").Value = "OK"

# Fila 38
$ws.Range("D38").Value = "12/17/2023"
$ws.Range("E38").Value = "SI"
$ws.Range("F38").Value = "-"
$ws.Range("G38").Value = "OK"

# Fila 39
$ws.Range("D39").Value = "12/17/2023"
$ws.Range("E39").Value = "SI"
$ws.Range("F39").Value = "-"
$ws.Range("G39").Value = "OK"

$ws.Range("H45").Select() | Out-Null
